# ADD: Battery voltage display
#
# Typography sheet: add a new "Medium" typography entry (row 8 of the
# Table7 typography table).
# Translation sheet: drop the per-locale "GB-TYPOGRAPHY" override column
# (now redundant since the base typography already says "Large"), and add
# two new translation rows - one more single-use id row, and the new
# "batteryVoltage" text.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Typography")
$ws2 = $wb.Worksheets.Item("Translation")

# ---------------------------------------------------------------------
# Typography sheet - new row 8: Medium / verdana.ttf / 12 / 4 / ? / (none)
# / 0x20-0x39 / (none)
# ---------------------------------------------------------------------
$ws1.Cells.Item(8, 2).Value = "Medium"
$ws1.Cells.Item(8, 3).Value = "verdana.ttf"
$ws1.Cells.Item(8, 4).Value = 12
$ws1.Cells.Item(8, 5).Value = 4
$ws1.Cells.Item(8, 6).Value = "?"
$ws1.Cells.Item(8, 8).Value = "0x20-0x39"

# ---------------------------------------------------------------------
# Translation sheet - the GB-TYPOGRAPHY (column G) override is removed;
# the base TYPOGRAPHY NAME column (C) now carries "Large" directly for
# the rows that used to rely on the override.
# ---------------------------------------------------------------------
$ws2.Cells.Item(3, 7).ClearContents()

$ws2.Cells.Item(4, 3).Value = "Large"
$ws2.Cells.Item(4, 7).ClearContents()

$ws2.Cells.Item(14, 3).Value = "Large"
$ws2.Cells.Item(14, 7).ClearContents()

# ---------------------------------------------------------------------
# Translation sheet - new row 15 (another single-use id placeholder row)
# ---------------------------------------------------------------------
$ws2.Cells.Item(15, 2).Value = "SingleUseId15"
$ws2.Cells.Item(15, 3).Value = "Small"
$ws2.Cells.Item(15, 4).Value = "Left"
$ws2.Cells.Item(15, 5).Value = "LTR"
# The wildcard value is the text "0" (matching the existing F12/F13
# cells), not the number 0, so copy an existing text "0" cell instead of
# typing a numeric-looking literal (which Excel would auto-convert to a
# Number).
$ws2.Range("F12").Copy($ws2.Range("F15"))
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Translation sheet - new row 16: batteryVoltage
# ---------------------------------------------------------------------
$ws2.Cells.Item(16, 2).Value = "batteryVoltage"
$ws2.Cells.Item(16, 3).Value = "Medium"
$ws2.Cells.Item(16, 4).Value = "Center"
$ws2.Cells.Item(16, 5).Value = "LTR"
$ws2.Cells.Item(16, 6).Value = "<text>V"
